$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally had a two-row, merged-cell header exported from pandas
# (group headers like "Tackles"/"Challenges"/"Blocks" merged across several
# columns, with "Unnamed: N_level_0" placeholders underneath). Clean this up
# by un-merging the header cells and writing a single flat header row that
# combines the group name with the sub-column name, while keeping the
# original two-level header available (but hidden) in row 2.

$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

$ws.Range("A1").Value = "Player ID"
$ws.Range("B1").Value = "Player"
$ws.Range("C1").Value = "#"
$ws.Range("D1").Value = "Nation"
$ws.Range("E1").Value = "Pos"
$ws.Range("F1").Value = "Age"
$ws.Range("G1").Value = "90s"
$ws.Range("H1").Value = "Tkl"
$ws.Range("I1").Value = "TklW"
$ws.Range("J1").Value = "Def 3rd"
$ws.Range("K1").Value = "Mid 3rd"
$ws.Range("L1").Value = "Att 3rd"
$ws.Range("M1").Value = "Cha"
$ws.Range("N1").Value = "Att"
$ws.Range("O1").Value = "Tkl%"
$ws.Range("P1").Value = "Lost"
$ws.Range("Q1").Value = "Blocks"
$ws.Range("R1").Value = "Sh"
$ws.Range("S1").Value = "Pass"
$ws.Range("T1").Value = "Int"
$ws.Range("U1").Value = "Tkl+Int"
$ws.Range("V1").Value = "Clr"
$ws.Range("W1").Value = "Err"

# Fill in the missing zeros that pandas/openpyxl had omitted for a few rows.
$ws.Range("O5").Value = 0
$ws.Range("O12").Value = 0
$ws.Range("O19").Value = 0

# Hide the now-redundant original sub-header row, the blank spacer row, and
# the trailing "16 Players" totals row, keeping the data intact but tidied.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true
$ws.Rows.Item(20).Hidden = $true

# Restore the last active selection that was in place when the cleaned file
# was saved.
[void]$ws.Range("O21").Select()
